$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 2).NumberFormat = "@"

$ws.Cells.Item(3, 1).Value = "Marc"
$ws.Cells.Item(3, 2).Value = "2026-02-12"
$ws.Cells.Item(3, 3).Value = "14:49:10"

$ws.Cells.Item(3, 2).Style = "Normal"
